$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet so it lands at the end
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "HypothyroidismElimination"

# ---- Header row ----
$ws.Range("A1").Value = "Recipe Id"
$ws.Range("B1").Value = "Recipe Name"
$ws.Range("C1").Value = "Recipe Category(Breakfast/lunch/snack/dinner)"
$ws.Range("D1").Value = "Food Category(Veg/non-veg/vegan/Jain)"
$ws.Range("E1").Value = "Ingredients"
$ws.Range("F1").Value = "Preparation Time"
$ws.Range("G1").Value = "Cooking Time"
$ws.Range("H1").Value = "Preparation method"
$ws.Range("I1").Value = "Nutrient values"
$ws.Range("J1").Value = "Targetted morbid conditions (Diabeties/Hypertension/Hypothyroidism)"
$ws.Range("K1").Value = "Recipe URL"

# ---- Data row for "Instant Rabri, Quick Rabdi" ----

# A2 ("2794") must be stored as text (like every other Recipe Id in this
# workbook), not auto-coerced to a number. Build it as a text-formula result
# first, then flatten to a plain value via copy/paste-values so the cell
# keeps the sheet's default style (no new number format is introduced).
$ws.Range("A2").Formula = "=""2794"""
$ws.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4163)

$ws.Range("B2").Value = "Instant Rabri, Quick Rabdi"
$ws.Range("C2").Value = "Lunch"
$ws.Range("D2").Value = "[Veg]"
$ws.Range("E2").Value = "For Instant Rabri`n3 cups full-fat milk`n2 fresh bread slices`n1 tbsp sugar`n1/2 cup condensed milk`nFor The Garnish`n1/4 tsp cardamom (elaichi) powder"
$ws.Range("F2").Value = "10 mins"
$ws.Range("G2").Value = "15 mins"
$ws.Range("H2").Value = "For instant rabri`nTo make instant rabri, remove the crusts of the bread slices and discard. Grind the bread slices in a food processor to make fresh bread crumbs and keep aside.`nBring the milk to boil in a broad non-stick pan.`nAdd the fresh bread crumbs, sugar and condensed milk, mix well and cook on a medium flame for approx. 8 to 10 minutes, while stirring continuously and scrapping the sides of the pan.`nAllow the instant rabri to cool completely and refrigerate for 2 to 3 hours.`nServe the instant rabri chilled garnished with cardamom powder."

# I2 (Nutrient values) is intentionally left blank for this recipe, but the
# cell still participates in the used range, so give it the default
# formatting (copied from the already-default-styled A1) without any value.
$ws.Range("A1").Copy()
$ws.Range("I2").PasteSpecial(-4122)

$ws.Range("J2").Value = "Diabetic"
$ws.Range("K2").Value = "https://tarladalal.com/instant-rabri-quick-rabdi-2794r"

$excel.CutCopyMode = 0

# Restore the original active sheet/selection (adding a sheet shouldn't
# change which tab is shown).
$wb.Worksheets.Item(1).Activate()

